$d = $word.ActiveDocument

# Remove the "Ben Jarman" and "Catherine Heard" author paragraphs entirely.
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    $text = $p.Range.Text.Trim()
    if ($text -eq "Ben Jarman" -or $text -eq "Catherine Heard") {
        $p.Range.Delete()
    }
}
